# Update cryptos list: apply new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so numeric-looking
# strings (e.g. "1.010", "335.32") are stored as text, not converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.999.09'
$ws.Range("D3").Value = '1.945.26'
$ws.Range("D4").Value = '1.010'
$ws.Range("D5").Value = '335.32'
$ws.Range("D7").Value = '0.4851'
$ws.Range("D8").Value = '0.4171'
$ws.Range("D9").Value = '0.08219'
$ws.Range("D10").Value = '1.019'
$ws.Range("D11").Value = '23.93'
$ws.Range("D12").Value = '1.960.50'
$ws.Range("D13").Value = '6.102'
$ws.Range("D14").Value = '7.346'
$ws.Range("D15").Value = '91.61'
$ws.Range("D16").Value = '0.06874'
$ws.Range("D19").Value = '17.90'
$ws.Range("D21").Value = '30.008.28'
$ws.Range("D22").Value = '5.663'
$ws.Range("D23").Value = '11.96'
$ws.Range("D24").Value = '2.192'
$ws.Range("D25").Value = '2.195.07'
$ws.Range("D26").Value = '6.599'
$ws.Range("D27").Value = '157.25'
$ws.Range("D28").Value = '20.14'
$ws.Range("D29").Value = '2.113'
$ws.Range("D30").Value = '121.41'
$ws.Range("D31").Value = '1.020'
$ws.Range("D32").Value = '0.09652'
$ws.Range("D33").Value = '5.638'
$ws.Range("D35").Value = '3.560'
$ws.Range("D36").Value = '0.06523'
$ws.Range("D37").Value = '0.02299'
$ws.Range("D38").Value = '1.222'
$ws.Range("D39").Value = '0.5985'
$ws.Range("D40").Value = '8.012'
$ws.Range("D41").Value = '10.76'
$ws.Range("D42").Value = '2.536'
$ws.Range("D43").Value = '0.1857'
$ws.Range("D44").Value = '1.277'
$ws.Range("D45").Value = '12.44'
$ws.Range("D46").Value = '0.07534'
$ws.Range("D47").Value = '0.5585'
$ws.Range("D48").Value = '1.989'
$ws.Range("D49").Value = '117.42'
$ws.Range("D50").Value = '2.437'
$ws.Range("D51").Value = '73.06'

# Restore the original (default) cell style on the Price column so no
# extraneous number-format styling is left applied to the cells.
$priceRange.Style = "Normal"

# Update the Volume(1h) column (plain text, already includes +/-/% and padding spaces)
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  +3.03%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("E9").Value = '  +0.64%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("E12").Value = '  +2.91%  '
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("E36").Value = '  +6.35%  '
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("E38").Value = '  +3.34%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  +6.00%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("E51").Value = '  +0.59%  '
